$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value = "2025-07-25T07:22:51+00:00"

# --- Elements sheet: repoint CodeSystem canonical URLs to the new MOS/NOS URLs ---
$wsElem = $wb.Worksheets.Item("Elements")

$wsElem.Range("Z3").Value  = "https://mos.esante.gouv.fr/NOS/TRE_R11-CiviliteExercice/FHIR/TRE-R11-CiviliteExercice?vs"
$wsElem.Range("Z8").Value  = "https://mos.esante.gouv.fr/NOS/TRE_R09-CategorieProfessionnelle/FHIR/TRE-R09-CategorieProfessionnelle?vs"
$wsElem.Range("Z11").Value = "https://mos.esante.gouv.fr/NOS/TRE_G09-DepartementOM/FHIR/TRE-G09-DepartementOM?vs"
$wsElem.Range("Z14").Value = "https://mos.esante.gouv.fr/NOS/TRE_R82-Ordre/FHIR/TRE-R82-Ordre?vs"
$wsElem.Range("Z17").Value = "https://mos.esante.gouv.fr/NOS/TRE_R33-StatutInscription/FHIR/TRE-R33-StatutInscription?vs"
$wsElem.Range("Z18").Value = "https://mos.esante.gouv.fr/NOS/TRE_G09-DepartementOM/FHIR/TRE-G09-DepartementOM?vs"
$wsElem.Range("Z20").Value = "https://mos.esante.gouv.fr/NOS/TRE_R03-AttributionParticuliere/FHIR/TRE-R03-AttributionParticuliere?vs"
$wsElem.Range("Z25").Value = "https://mos.esante.gouv.fr/NOS/TRE_R223-NatCycleForm/FHIR/TRE-R223-NatCycleForm?vs"
$wsElem.Range("Z26").Value = "https://mos.esante.gouv.fr/NOS/TRE_R224-NiveauFormAcquis/FHIR/TRE-R224-NiveauFormAcquis?vs"
$wsElem.Range("Z27").Value = "https://mos.esante.gouv.fr/NOS/TRE_R225-AnneeUniversitaire/FHIR/TRE-R225-AnneeUniversitaire?vs"

# The longer replacement URLs widen the "Binding Value Set" column; mirror that
# by growing column Z's stored width to match the new best-fit content width.
$wsElem.Columns.Item(26).ColumnWidth = 86.5
